$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns stay as text, matching the
# original inline-string cell types (so e.g. "593.45" is not coerced to a number).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '67.896.01'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '3.777.32'
$ws.Range('E3').Value = '  -0.32%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '593.45'
$ws.Range('E5').Value = '  -0.62%  '
$ws.Range('D6').Value = '166.58'
$ws.Range('E6').Value = '  -1.39%  '
$ws.Range('D7').Value = '3.774.26'
$ws.Range('E7').Value = '  -0.44%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '0.519'
$ws.Range('E9').Value = '  -0.89%  '
$ws.Range('D10').Value = '0.159'
$ws.Range('E10').Value = '  -1.42%  '
$ws.Range('D11').Value = '6.38'
$ws.Range('E11').Value = '  -2.13%  '
$ws.Range('D12').Value = '0.451'
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('D13').Value = '0.0000257'
$ws.Range('E13').Value = '  -3.53%  '
$ws.Range('D14').Value = '35.98'
$ws.Range('E14').Value = '  -2.38%  '
$ws.Range('D15').Value = '4.416.35'
$ws.Range('E15').Value = '  -0.16%  '
$ws.Range('D16').Value = '3.795.53'
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('D17').Value = '68.008.88'
$ws.Range('E17').Value = '  -0.91%  '
$ws.Range('E18').Value = '  -3.30%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').Value = '0.112'
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = '6.94'
$ws.Range('E20').Value = '  -2.08%  '
$ws.Range('D21').Value = '10.53'
$ws.Range('E21').Value = '  -3.28%  '
$ws.Range('D22').Value = '461.09'
$ws.Range('E22').Value = '  -2.34%  '
$ws.Range('D23').Value = '0.695'
$ws.Range('E23').Value = '  -1.66%  '
$ws.Range('D24').Value = '0.0000152'
$ws.Range('E24').Value = '  +4.60%  '
$ws.Range('D25').Value = '83.47'
$ws.Range('E25').Value = '  -1.51%  '
$ws.Range('D26').Value = '2.14'
$ws.Range('E26').Value = '  -5.14%  '
$ws.Range('D27').Value = '11.84'
$ws.Range('E27').Value = '  -3.72%  '
$ws.Range('D28').Value = '10.03'
$ws.Range('E28').Value = '  -2.25%  '
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('D30').Value = '2.75'
$ws.Range('E30').Value = '  -2.03%  '
$ws.Range('D31').Value = '29.80'
$ws.Range('E31').Value = '  -1.66%  '
$ws.Range('D32').Value = '7.18'
$ws.Range('E32').Value = '  -4.48%  '
$ws.Range('D33').Value = '2.15'
$ws.Range('E33').Value = '  -4.35%  '
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +0.43%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').Value = '9.06'
$ws.Range('E35').Value = '  -2.30%  '
$ws.Range('D36').Value = '3.730.32'
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('D37').Value = '0.1000'
$ws.Range('E37').Value = '  -2.56%  '
$ws.Range('D38').Value = '3.43'
$ws.Range('E38').Value = '  -1.86%  '
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('D40').Value = '0.992'
$ws.Range('E40').Value = '  -1.30%  '
$ws.Range('D41').Value = '5.74'
$ws.Range('E41').Value = '  -2.01%  '
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('D44').Value = '0.299'
$ws.Range('E44').Value = '  -3.83%  '
$ws.Range('D45').Value = '43.16'
$ws.Range('E45').Value = '  +3.20%  '
$ws.Range('D46').Value = '46.83'
$ws.Range('E46').Value = '  +2.90%  '
$ws.Range('D47').Value = '1.90'
$ws.Range('E47').Value = '  -4.10%  '
$ws.Range('D48').Value = '8.34'
$ws.Range('E48').Value = '  -3.45%  '
$ws.Range('D49').Value = '146.84'
$ws.Range('D50').Value = '386.54'
$ws.Range('E50').Value = '  -4.76%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.757.88'
$ws.Range('E51').Value = '  +3.32%  '
